# Add a "2022-Q1" per-fund holdings sheet right before the "总计" (Total)
# summary sheet, and prepend a matching 2022-Q1 summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right before "总计".
#    NOTE: the sheet object passed as the "Before" argument to Add() gets
#    rebound by this host to the newly inserted sheet, so never reuse that
#    variable for anything else afterwards -- re-fetch "总计" fresh instead.
# ---------------------------------------------------------------------------
$beforeAnchor = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($beforeAnchor)
$newSheet.Name = "2022-Q1"

# Copy the header/index-column formatting (bold, centered, bordered -- the
# style already used across all the per-quarter detail sheets) instead of
# inventing a new style.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A17").PasteSpecial(-4122)

# Header row.
$newHeaders = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $newHeaders.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $newHeaders[$c]
}

# Per-fund holdings detail for 2022-Q1 (code, name, fund size, stock
# position, position share, held market value (100M yuan), position rank).
$fundRows = @(
    @("008515", "富兰克林国海基本面优选混合", "13.82", "89.34", "7.02", "0.9702", 2),
    @("513690", "博时恒生港股通高股息率ETF", "4.60", "99.64", "2.33", "0.1072", 10),
    @("005051", "上投摩根标普港股通低波红利指数A", "4.02", "92.23", "2.55", "0.1025", 6),
    @("005052", "上投摩根标普港股通低波红利指数C", "2.61", "92.23", "2.55", "0.0666", 6),
    @("005576", "华泰柏瑞新金融地产灵活配置混合", "0.79", "94.50", "3.40", "0.0269", 8),
    @("007751", "景顺长城中证沪港深红利成长低波动指数A", "0.83", "91.29", "2.81", "0.0233", 5),
    @("004316", "前海开源沪港深裕鑫灵活配置混合A", "0.64", "90.55", "3.05", "0.0195", 10),
    @("004317", "前海开源沪港深裕鑫灵活配置混合C", "0.47", "90.55", "3.05", "0.0143", 10),
    @("004098", "前海开源港股通股息率50强股票", "0.34", "88.92", "3.86", "0.0131", 3),
    @("005702", "恒生前海港股通高股息低波动指数", "0.29", "94.14", "2.38", "0.0069", 8),
    @("006658", "财通中证香港红利等权投资指数A", "0.20", "90.59", "3.41", "0.0068", 7),
    @("001824", "博时沪港深成长企业混合", "0.15", "93.41", "3.36", "0.0050", 8),
    @("501307", "银河中证沪港深高股息指数（LOF）A", "0.19", "91.35", "1.42", "0.0027", 10),
    @("006659", "财通中证香港红利等权投资指数C", "0.05", "90.59", "3.41", "0.0017", 7),
    @("007760", "景顺长城中证沪港深红利成长低波动指数C", "0.06", "91.29", "2.81", "0.0017", 5),
    @("501308", "银河中证沪港深高股息指数（LOF）C", "0.01", "91.35", "1.42", "0.0001", 10)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Prepend a 2022-Q1 summary row to the "总计" sheet, pushing the other
#    rows down and renumbering the index column (A) to match.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$origLastRow = $totalSheet.UsedRange.Rows.Count

$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

for ($r = 3; $r -le ($origLastRow + 1); $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 16
$totalSheet.Cells.Item(2, 4).Value = 1.37
